$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Coin / Link text cells (B, C columns) ---
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"

# --- Update Price cells (D column) as TEXT, avoiding numeric auto-conversion ---
$dCells = @("D2","D3","D4","D6","D10","D11","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D25","D26","D27","D29","D30","D31","D32","D36","D37","D41","D46","D47","D48","D49","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("D2").Value = "43.242.38"
$ws.Range("D3").Value = "2.599.13"
$ws.Range("D4").Value = "0.999"
$ws.Range("D6").Value = "98.29"
$ws.Range("D10").Value = "35.99"
$ws.Range("D11").Value = "0.0815"
$ws.Range("D13").Value = "2.989.58"
$ws.Range("D14").Value = "0.109"
$ws.Range("D15").Value = "2.579.50"
$ws.Range("D16").Value = "15.29"
$ws.Range("D17").Value = "0.851"
$ws.Range("D18").Value = "43.304.31"
$ws.Range("D20").Value = "12.81"
$ws.Range("D21").Value = "0.0₃0970"
$ws.Range("D22").Value = "69.70"
$ws.Range("D23").Value = "254.93"
$ws.Range("D25").Value = "2.10"
$ws.Range("D26").Value = "27.31"
$ws.Range("D27").Value = "1.00"
$ws.Range("D29").Value = "41.12"
$ws.Range("D30").Value = "10.34"
$ws.Range("D31").Value = "5.89"
$ws.Range("D32").Value = "156.34"
$ws.Range("D36").Value = "2.70"
$ws.Range("D37").Value = "18.88"
$ws.Range("D41").Value = "22.90"
$ws.Range("D46").Value = "2.015.04"
$ws.Range("D47").Value = "9.00"
$ws.Range("D48").Value = "2.845.13"
$ws.Range("D49").Value = "83.64"
$ws.Range("D51").Value = "74.97"
foreach ($addr in $dCells) { $ws.Range($addr).Style = "Normal" }

# --- Update Volume(1h) percentage cells (E column) ---
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  +4.14%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("E25").Value = "  +4.34%  "
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  +6.21%  "
$ws.Range("E34").Value = "  +2.94%  "
$ws.Range("E35").Value = "  +3.79%  "
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("E39").Value = "  +10.03%  "
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("E42").Value = "  +6.62%  "
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("E50").Value = "  +4.70%  "
$ws.Range("E51").Value = "  +2.03%  "
